$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.468200000000003
$ws.Range("A3").Value = -21.8729
$ws.Range("C3").Value = -11.1143
$ws.Range("C12").Value = -11.61229999999999
$ws.Range("A14").Value = -21.80090000000001
$ws.Range("A16").Value = -21.974
$ws.Range("B18").Value = 6.995899999999997
$ws.Range("A21").Value = -20.12069999999998
$ws.Range("A23").Value = -20.73339999999997
$ws.Range("B24").Value = 6.634500000000002
$ws.Range("C24").Value = -12.50599999999999
$ws.Range("A25").Value = -21.90909999999999
$ws.Range("B25").Value = 5.942999999999999
$ws.Range("C25").Value = -13.52879999999998
$ws.Range("A26").Value = -20.99659999999996
$ws.Range("B27").Value = 5.806800000000003
$ws.Range("A29").Value = -20.75499999999997
$ws.Range("B30").Value = 5.647000000000003
$ws.Range("B31").Value = 5.5886
$ws.Range("B39").Value = 9.545599999999999
$ws.Range("A40").Value = -19.64350000000001
$ws.Range("C41").Value = -12.2182
$ws.Range("B42").Value = 10.183
$ws.Range("B48").Value = 4.874700000000003
$ws.Range("C50").Value = -13.48349999999999
$ws.Range("B51").Value = 5.576100000000002
$ws.Range("B52").Value = 5.475900000000001
$ws.Range("A53").Value = -21.3639
$ws.Range("C53").Value = -10.2413
$ws.Range("B55").Value = 6.634199999999993
$ws.Range("B56").Value = 5.447700000000001
$ws.Range("C56").Value = -12.30289999999999
$ws.Range("A57").Value = -21.96079999999998
$ws.Range("B57").Value = 5.343799999999999
$ws.Range("C57").Value = -12.86169999999999
$ws.Range("C58").Value = -13.3912
$ws.Range("A59").Value = -22.181
$ws.Range("B60").Value = 5.8217
$ws.Range("C61").Value = -13.05589999999999
$ws.Range("C63").Value = -11.96129999999999
$ws.Range("C64").Value = -11.67649999999999
$ws.Range("A65").Value = -21.76639999999999
$ws.Range("A69").Value = -21.5631
$ws.Range("C70").Value = -11.2118
$ws.Range("C72").Value = -11.7692
$ws.Range("B73").Value = 8.785799999999995
$ws.Range("B74").Value = 9.544199999999991
$ws.Range("A79").Value = -20.6155
$ws.Range("A83").Value = -21.75339999999999
$ws.Range("C86").Value = -13.2529
$ws.Range("B89").Value = 5.016299999999995
$ws.Range("C89").Value = -10.12590000000001
$ws.Range("B90").Value = 5.6294
$ws.Range("A91").Value = -21.42500000000001
$ws.Range("B92").Value = 5.082299999999991
$ws.Range("A93").Value = -20.83999999999998
$ws.Range("C98").Value = -11.619
$ws.Range("A100").Value = -21.8063
$ws.Range("C100").Value = -13.18609999999998
$ws.Range("C102").Value = -12.4171
